$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 279. This shifts all existing rows (279..291)
# down to (282..294), automatically producing the trailing duplicated block
# that appears in the target file (rows 292-294, the old Winter Nelis block).
$ws.Rows("279:281").Insert()

# Row 279: new weekly entry - Packham's Triumph / Especial
$ws.Cells.Item(279, 1).Value = 8
$ws.Cells.Item(279, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(279, 3).Value = "Coquimbo"
$ws.Cells.Item(279, 4).Value = 44448
$ws.Cells.Item(279, 5).Value = 4
$ws.Cells.Item(279, 6).Value = "Fruta"
$ws.Cells.Item(279, 7).Value = 100104
$ws.Cells.Item(279, 8).Value = "Frutos de pepita"
$ws.Cells.Item(279, 9).Value = 100104005
$ws.Cells.Item(279, 10).Value = "Pera"
$ws.Cells.Item(279, 11).Value = "Packham's Triumph"
$ws.Cells.Item(279, 12).Value = "Especial"
$ws.Cells.Item(279, 13).Value = 16
$ws.Cells.Item(279, 14).Value = 255000
$ws.Cells.Item(279, 15).Value = 260000
$ws.Cells.Item(279, 16).Value = 257500
$ws.Cells.Item(279, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(279, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(279, 19).Value = 572
$ws.Cells.Item(279, 20).Value = 450

# Row 280: new weekly entry - Packham's Triumph / Primera
$ws.Cells.Item(280, 1).Value = 8
$ws.Cells.Item(280, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(280, 3).Value = "Coquimbo"
$ws.Cells.Item(280, 4).Value = 44448
$ws.Cells.Item(280, 5).Value = 4
$ws.Cells.Item(280, 6).Value = "Fruta"
$ws.Cells.Item(280, 7).Value = 100104
$ws.Cells.Item(280, 8).Value = "Frutos de pepita"
$ws.Cells.Item(280, 9).Value = 100104005
$ws.Cells.Item(280, 10).Value = "Pera"
$ws.Cells.Item(280, 11).Value = "Packham's Triumph"
$ws.Cells.Item(280, 12).Value = "Primera"
$ws.Cells.Item(280, 13).Value = 20
$ws.Cells.Item(280, 14).Value = 235000
$ws.Cells.Item(280, 15).Value = 240000
$ws.Cells.Item(280, 16).Value = 237500
$ws.Cells.Item(280, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(280, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(280, 19).Value = 528
$ws.Cells.Item(280, 20).Value = 450

# Row 281: new weekly entry - Packham's Triumph / Segunda
$ws.Cells.Item(281, 1).Value = 8
$ws.Cells.Item(281, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(281, 3).Value = "Coquimbo"
$ws.Cells.Item(281, 4).Value = 44448
$ws.Cells.Item(281, 5).Value = 4
$ws.Cells.Item(281, 6).Value = "Fruta"
$ws.Cells.Item(281, 7).Value = 100104
$ws.Cells.Item(281, 8).Value = "Frutos de pepita"
$ws.Cells.Item(281, 9).Value = 100104005
$ws.Cells.Item(281, 10).Value = "Pera"
$ws.Cells.Item(281, 11).Value = "Packham's Triumph"
$ws.Cells.Item(281, 12).Value = "Segunda"
$ws.Cells.Item(281, 13).Value = 20
$ws.Cells.Item(281, 14).Value = 205000
$ws.Cells.Item(281, 15).Value = 210000
$ws.Cells.Item(281, 16).Value = 207500
$ws.Cells.Item(281, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(281, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(281, 19).Value = 461
$ws.Cells.Item(281, 20).Value = 450

# Ensure date column (D) is formatted as date like the rest of column D
$ws.Range("D279:D281").NumberFormat = "YYYY-MM-DD HH:MM:SS"
